$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted at row 180 (pushing the
# existing rows 180-261 down to 181-262, which keeps every other row's
# data unchanged other than the shift).
$ws.Rows("180:180").Insert()

# Fill in the new row with the latest reading.
$ws.Range("A180").Value = 8
$ws.Range("B180").Value = "Terminal La Palmera de La Serena"
$ws.Range("C180").Value = "Coquimbo"
$ws.Range("D180").Value = 45202
$ws.Range("E180").Value = 4
$ws.Range("F180").Value = 100112040
$ws.Range("G180").Value = "Cilantro"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 2400
$ws.Range("K180").Value = 1500
$ws.Range("L180").Value = 2000
$ws.Range("M180").Value = 1750
$ws.Range("N180").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O180").Value = "Provincia del Elqu" + [char]0x00ED
$ws.Range("P180").Value = 1167
$ws.Range("Q180").Value = 1.5
$ws.Range("R180").Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format as the
# rest of column D.
$ws.Range("D180").NumberFormat = $ws.Range("D181").NumberFormat
